$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 600, pushing existing rows 600-650 down to 601-651
$ws.Rows("600").Insert()

# Populate the newly inserted row 600 with data
$ws.Range("A600").Value = 4
$ws.Range("B600").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C600").Value = "Los Lagos"
$ws.Range("D600").Value = 44931
$ws.Range("E600").Value = 10
$ws.Range("F600").Value = "Fruta"
$ws.Range("G600").Value = 100106
$ws.Range("H600").Value = "Oleaginosos"
$ws.Range("I600").Value = 100106002
$ws.Range("J600").Value = "Palta"
$ws.Range("K600").Value = "Hass"
$ws.Range("L600").Value = "Primera"
$ws.Range("M600").Value = 100
$ws.Range("N600").Value = 4500
$ws.Range("O600").Value = 4500
$ws.Range("P600").Value = 4500
$ws.Range("Q600").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R600").Value = "Provincia de Quillota"
$ws.Range("S600").Value = 4500
$ws.Range("T600").Value = 1
